# Apply the "Force" row insertion + related refactor to the Server.xlsx sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row above the existing "Cache" row (row 8), shifting
#    everything below it (header row, GameServer..Tutorial data rows) down by one.
$ws.Rows("8:8").Insert()

# 2. The new row 8 should look like its neighbours (row 7 "Save" / the old
#    row 8 "Cache") - same fill/border/font/number format - so copy formats
#    from row 7 into the freshly inserted row 8.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

# 3. Give the new row its label and boolean-list formulas (same pattern as
#    every other row in this block: FALSE() across B:I).
$ws.Range("A8").Value = "Force"
$ws.Range("B8:I9").Formula = "=FALSE()"

# 4. The header row (now row 10) grew a touch taller after the edit.
$ws.Rows("10:10").RowHeight = 41.65

# 5. Re-anchor the frozen pane on the row below the header (now row 11) and
#    leave the active selection on the new "Force" row's label cell.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A11").Select()
$win.FreezePanes = $true
$ws.Range("A9").Select()

# 6. Cosmetic: the default cell style carries the English "Normal" name in
#    the refreshed workbook rather than the localized "常规".
$style = $wb.Styles.Item(1)
$style.Name = "Normal"

Write-Output "edit applied"
